$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.432.24'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.84%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.980.10'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.72%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '502.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.75'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +9.21%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +7.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.33'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +10.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.106'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +11.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.350'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.21%  '
$ws.Range('E12').Value = '  +2.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.488.68'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.13'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +12.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '56.421.93'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.76%  '
$ws.Range('E16').Value = '  +14.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.978.26'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.66'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +11.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.34'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.77'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +11.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.53'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.83%  '
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('E23').Value = '  +5.95%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.91'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.24%  '
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.163'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0₃0891'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +12.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.52'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +11.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.79'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +13.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.19'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.29%  '
$ws.Range('E32').Value = '  +10.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.43'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '156.97'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +8.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.45'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.61%  '
$ws.Range('E36').Value = '  +5.74%  '
$ws.Range('E37').Value = '  +4.12%  '
$ws.Range('E38').Value = '  +11.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '22.98'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.013.18'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.92%  '
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.21'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.14%  '
$ws.Range('E43').Value = '  +6.88%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.247.56'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.991'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.85%  '
$ws.Range('E46').Value = '  +8.09%  '
$ws.Range('E47').Value = '  +5.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.94'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +25.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0234'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +12.67%  '
$ws.Range('E50').Value = '  +9.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.97'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.96%  '
